$d = $word.ActiveDocument

# Locate the paragraph whose text is "  d = a" - the body of the
# "if (a < b) { ... }" block - so we can add the missing "else" clause
# right after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^\s*d = a\s*$") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the '  d = a' paragraph"
}

# Insert a new paragraph right after "  d = a" and give it the closing
# brace for the "if" block.
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pClose = $target.Next()
$pClose.Range.InsertAfter("}")

# Insert another new paragraph after that one for the "else {" line.
$pClose.Range.InsertParagraphAfter()
$pElse = $pClose.Next()
$pElse.Range.InsertAfter("else {")

# Insert a final new paragraph for the else-branch body, "  d = b".
$pElse.Range.InsertParagraphAfter()
$pBody = $pElse.Next()
$pBody.Range.InsertAfter("  d")
$pBody.Range.InsertAfter(" ")
$pBody.Range.InsertAfter("=")
$pBody.Range.InsertAfter(" ")
$pBody.Range.InsertAfter("b")

Write-Output "Inserted else clause; paragraph count now $($d.Paragraphs.Count)"
